# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$nl = [char]10
$newText = "Conversión del día 💰" + $nl + `
  "✅ Dólar paralelo: 68" + $nl + `
  "" + $nl + `
  "Binance" + $nl + `
  "✅ 1000 Bs = 14.36 = 59382.63 pesos" + $nl + `
  "✅ 59382.63 pesos = 14.24 = 960.76 Bs" + $nl + `
  "" + $nl + `
  "Promedio competencia" + $nl + `
  "✅ Tasa pesos: 20" + $nl + `
  "✅ Tasa Bs: 20" + $nl + `
  "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate cells on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 69.65000000000001
$wsTasas.Range("O10").Value = 4136

$wsTasas.Range("N12").Value = 4169
$wsTasas.Range("O12").Value = 67.45099999999999
